$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the config path value from E2 (keep its style/format)
$ws.Range("E2").ClearContents()

# Move the active selection to D10 as in the edited workbook
$ws.Range("D10").Select()
